# Applies the "fINAL COMMIT FOR LOGIN,LOGINOUT,GET ALL PROGRAM/ID/USER" edit
# to the Program sheet of the LMS hackathon test-data workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# --- GET_ALLUSERS_PROGRAM_21 row (22): Endpoint now points at the invalid variant ---
$ws.Range("F22").Value = "/allProgramsWithUsersInvalid"

# --- GET_ALL_PROGRAM_TC11 row (12): Action column now carries the schema-validation tag ---
$ws.Range("G12").Value = "validateSchemaProgram"

# --- GET_ALL_PROGRAM_14 row (16): Action column switched from AddProgramId -> validateSchemaProgram
#     (this also adopts the bordered "text" style used elsewhere in the block, i.e. style 6) ---
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Font.Size = 12
$ws.Range("G16").Value = "validateSchemaProgram"

# --- GET_ID_PROGRAM_17/18 rows (18/19): Action column text simplified ---
$ws.Range("G18").Value = "InvalidUri"
$ws.Range("G19").Value = "NoAuth"

# --- GET_ALLUSERS_PROGRAM_20 row (21): Action column now carries the schema-validation tag
#     (also adopts style 6, same as G12/G16) ---
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Font.Size = 12
$ws.Range("G21").Value = "validateSchemaProgram"

# --- Row 2: programName / programDescription values got a "TCNINE" suffix ---
$ws.Range("B2").Value = "WellsTCNINE"
$ws.Range("D2").Value = "testerTCNINE"

# Reflect the cursor/selection that was left on the sheet after the edits.
$ws.Activate() | Out-Null
$ws.Range("G21").Select() | Out-Null
